# Update test data and add one more test case.
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the TestCases and Sheet3 worksheets, leaving only TestData.
$wb.Worksheets.Item("TestCases").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

$ws = $wb.Worksheets.Item("TestData")

# Refresh the test data table: two columns (Product_Category / Search_Text)
# with a single data row (Books / Selenium).
$ws.Range("A2").Value = "Product_Category"
$ws.Range("B2").Value = "Search_Text"
$ws.Range("C2").Clear()

$ws.Range("A3").Value = "Books"
$ws.Range("B3").Value = "Selenium"
$ws.Range("C3").Clear()

# Move the selection as recorded by the last save.
$ws.Range("B11").Select() | Out-Null
